# Update LR-pair statistics following Dr Hou advice (natmi re-run)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 8.200698666666668
$ws.Range("H2").Value = 24.602096
$ws.Range("I2").Value = 0.3423472217473603
$ws.Range("J2").Value = 0.3423472217473603
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 78.370804
$ws.Range("N2").Value = 235.112412
$ws.Range("O2").Value = 0.9256182775132763
$ws.Range("P2").Value = 0.9256182775132761
$ws.Range("Q2").Value = 642.6953478683948
$ws.Range("R2").Value = 5784.258130815553
$ws.Range("S2").Value = 0.3168828457052473
$ws.Range("T2").Value = 0.3168828457052472

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 8.200698666666668
$ws.Range("H3").Value = 24.602096
$ws.Range("I3").Value = 0.3423472217473603
$ws.Range("J3").Value = 0.3423472217473603
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.309992333333333
$ws.Range("N3").Value = 6.929977
$ws.Range("O3").Value = 0.02728275091638557
$ws.Range("P3").Value = 0.02728275091638557
$ws.Range("Q3").Value = 18.94355104797689
$ws.Range("R3").Value = 170.491959431792
$ws.Range("S3").Value = 0.009340173977849846
$ws.Range("T3").Value = 0.009340173977849848

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 8.200698666666668
$ws.Range("H4").Value = 24.602096
$ws.Range("I4").Value = 0.3423472217473603
$ws.Range("J4").Value = 0.3423472217473603
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.987804
$ws.Range("N4").Value = 11.963412
$ws.Range("O4").Value = 0.04709897157033827
$ws.Range("P4").Value = 0.04709897157033827
$ws.Range("Q4").Value = 32.70277894572801
$ws.Range("R4").Value = 294.325010511552
$ws.Range("S4").Value = 0.01612420206426321
$ws.Range("T4").Value = 0.01612420206426321

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 10.62814333333333
$ws.Range("H5").Value = 31.88443
$ws.Range("I5").Value = 0.4436835799477486
$ws.Range("J5").Value = 0.4436835799477487
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 78.370804
$ws.Range("N5").Value = 235.112412
$ws.Range("O5").Value = 0.9256182775132763
$ws.Range("P5").Value = 0.9256182775132761
$ws.Range("Q5").Value = 832.9361380605734
$ws.Range("R5").Value = 7496.42524254516
$ws.Range("S5").Value = 0.410681631032159
$ws.Range("T5").Value = 0.410681631032159

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 10.62814333333333
$ws.Range("H6").Value = 31.88443
$ws.Range("I6").Value = 0.4436835799477486
$ws.Range("J6").Value = 0.4436835799477487
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.309992333333333
$ws.Range("N6").Value = 6.929977
$ws.Range("O6").Value = 0.02728275091638557
$ws.Range("P6").Value = 0.02728275091638557
$ws.Range("Q6").Value = 24.55092961756778
$ws.Range("R6").Value = 220.95836655811
$ws.Range("S6").Value = 0.01210490859740467
$ws.Range("T6").Value = 0.01210490859740467

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 10.62814333333333
$ws.Range("H7").Value = 31.88443
$ws.Range("I7").Value = 0.4436835799477486
$ws.Range("J7").Value = 0.4436835799477487
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.987804
$ws.Range("N7").Value = 11.963412
$ws.Range("O7").Value = 0.04709897157033827
$ws.Range("P7").Value = 0.04709897157033827
$ws.Range("Q7").Value = 42.38295249724001
$ws.Range("R7").Value = 381.44657247516
$ws.Range("S7").Value = 0.02089704031818492
$ws.Range("T7").Value = 0.02089704031818492

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 5.125489
$ws.Range("H8").Value = 15.376467
$ws.Range("I8").Value = 0.2139691983048911
$ws.Range("J8").Value = 0.2139691983048911
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 78.370804
$ws.Range("N8").Value = 235.112412
$ws.Range("O8").Value = 0.9256182775132763
$ws.Range("P8").Value = 0.9256182775132761
$ws.Range("Q8").Value = 401.688693823156
$ws.Range("R8").Value = 3615.198244408404
$ws.Range("S8").Value = 0.1980538007758699
$ws.Range("T8").Value = 0.1980538007758699

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 5.125489
$ws.Range("H9").Value = 15.376467
$ws.Range("I9").Value = 0.2139691983048911
$ws.Range("J9").Value = 0.2139691983048911
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.309992333333333
$ws.Range("N9").Value = 6.929977
$ws.Range("O9").Value = 0.02728275091638557
$ws.Range("P9").Value = 0.02728275091638557
$ws.Range("Q9").Value = 11.83984029458433
$ws.Range("R9").Value = 106.558562651259
$ws.Range("S9").Value = 0.005837668341131051
$ws.Range("T9").Value = 0.005837668341131052

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 5.125489
$ws.Range("H10").Value = 15.376467
$ws.Range("I10").Value = 0.2139691983048911
$ws.Range("J10").Value = 0.2139691983048911
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.987804
$ws.Range("N10").Value = 11.963412
$ws.Range("O10").Value = 0.04709897157033827
$ws.Range("P10").Value = 0.04709897157033827
$ws.Range("Q10").Value = 20.439445536156
$ws.Range("R10").Value = 183.955009825404
$ws.Range("S10").Value = 0.01007772918789014
$ws.Range("T10").Value = 0.01007772918789014
